$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65

# Column A holds the purchase date as literal text (e.g. "01/04/2026" in the
# row above), not a real date serial. Typing a date-shaped string straight
# into a General-formatted cell makes Excel auto-convert it to a date, so we
# briefly mark the cell as Text first, enter the value, then restore the
# "Normal" cell style so no stray number-format survives on the new cell.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "01/11/2026"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = 0.0005434200000000028
$ws.Range("C$row").Value = 91089.76482278854
$ws.Range("D$row").Value = 50
